# Apply the "cleaned up clean rating file" update:
#  - correct a handful of mis-keyed Moody's/S&P/Fitch rating cells
#  - add a new country row (Algeria) at the bottom of the table
#  - update the active cell selection to match the saved workbook state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix mismatched rating cells -------------------------------------------------
$ws.Range("C18").Value  = "Baa1"   # Botswana - moodys
$ws.Range("D29").Value  = "BB"     # Costa Rica - s&p
$ws.Range("D38").Value  = "B"      # Egypt - s&p
$ws.Range("D44").Value  = "A+"     # France - s&p
$ws.Range("C49").Value  = "Caa1"   # Ghana - moodys
$ws.Range("E52").Value  = "BB+"    # Guatemala - fitch
$ws.Range("E71").Value  = "CCC+"   # Lao P.D.R. - fitch
$ws.Range("D78").Value  = "B- *-"  # Madagascar - s&p (new rating string)
$ws.Range("D84").Value  = "BB-"    # Moldova - s&p
$ws.Range("C85").Value  = "B1"     # Mongolia - moodys
$ws.Range("C111").Value = "Caa1"   # Senegal - moodys
$ws.Range("E116").Value = "A+"     # Slovenia - fitch

# --- Append the new row for Algeria ----------------------------------------------
$ws.Range("A140").Value = 612
$ws.Range("B140").Value = "Algeria"
$ws.Range("C140").Value = "NR"
$ws.Range("D140").Value = "NR"
$ws.Range("E140").Value = "NR"

# --- Update the active selection to match the saved workbook state ---------------
[void]$ws.Range("G16").Select()
